# TimeInPhase.xlsx — "Corregí los psp's de la ClaseConfiguración"
#
# The sheet/connection/query were all originally called "excel". This edit
# renames the sheet to "excel(1)" (and, as a consequence, the defined name's
# RefersTo formula picks up the quoted sheet reference automatically), gives
# the local-scope defined name a distinct identifier ("excel_1", since Excel
# does not allow a defined name to collide with a sheet name), refreshes the
# "Reporte generado..." timestamp string pulled from the web query, and
# relabels the underlying web query / connection to match the "(1)" suffix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet: excel -> excel(1)
#    Excel automatically re-quotes the sheet reference inside any defined
#    name's RefersTo formula that pointed at it (excel!... -> 'excel(1)'!...).
$ws.Name = "excel(1)"

# 2) Rename the local defined name: excel -> excel_1
#    (it now refers to 'excel(1)'!$A$1:$E$27, already updated by the rename above)
foreach ($n in $wb.Names) {
    if ($n.Name -eq "excel(1)!excel" -or $n.Name -eq "excel") {
        $n.Name = "excel_1"
    }
}

# 3) Refresh the generated-report timestamp text (result of the last web query refresh)
$ws.Range("A25").Value = "Reporte generado a las 11:22 AM el 5/12/2018"

# 4) Best-effort: relabel the legacy web query connection + its query table to
#    match the new "(1)" suffixed name (mirrors what Excel does when you open
#    a second copy of the same query/connection in one workbook).
try {
    foreach ($conn in $wb.Connections) {
        if ($conn.Name -eq "excel") {
            $conn.Name = "excel(1)"
            $conn.ODCFile = "C:\Users\AxeelZR\Downloads\excel(1).iqy"
        }
    }
} catch {
}

try {
    foreach ($qt in $ws.QueryTables) {
        if ($qt.Name -eq "excel") {
            $qt.Name = "excel(1)"
        }
    }
} catch {
}
